$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact literal text value (avoids Excel auto-numeric
# coercion / float rounding for numeric-looking strings like "581.19").
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" "69.380.89"

# Row 3
Set-TextCell $ws "D3" "3.539.43"
Set-TextCell $ws "E3" "  -2.82%  "

# Row 4
Set-TextCell $ws "E4" "  -0.08%  "

# Row 5
Set-TextCell $ws "D5" "581.19"
Set-TextCell $ws "E5" "  +0.69%  "

# Row 6
Set-TextCell $ws "D6" "172.65"
Set-TextCell $ws "E6" "  -1.98%  "

# Row 7
Set-TextCell $ws "B7" "XRP"
Set-TextCell $ws "C7" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws "D7" "0.612"
Set-TextCell $ws "E7" "  +0.07%  "

# Row 8
Set-TextCell $ws "B8" "LidoStakedEther"
Set-TextCell $ws "C8" "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextCell $ws "D8" "3.530.94"
Set-TextCell $ws "E8" "  -2.85%  "

# Row 9
Set-TextCell $ws "E9" "  -0.01%  "

# Row 10
Set-TextCell $ws "D10" "0.190"
Set-TextCell $ws "E10" "  -3.32%  "

# Row 11
Set-TextCell $ws "D11" "6.77"
Set-TextCell $ws "E11" "  -1.37%  "

# Row 12
Set-TextCell $ws "D12" "0.584"
Set-TextCell $ws "E12" "  -3.15%  "

# Row 13
Set-TextCell $ws "D13" "47.43"
Set-TextCell $ws "E13" "  -2.21%  "

# Row 14
Set-TextCell $ws "E14" "  -4.43%  "

# Row 15
Set-TextCell $ws "D15" "4.109.81"
Set-TextCell $ws "E15" "  -2.79%  "

# Row 16
Set-TextCell $ws "D16" "8.54"
Set-TextCell $ws "E16" "  -3.72%  "

# Row 17
Set-TextCell $ws "D17" "629.43"
Set-TextCell $ws "E17" "  -5.85%  "

# Row 18
Set-TextCell $ws "B18" "WrappedEther"
Set-TextCell $ws "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D18" "3.553.19"
Set-TextCell $ws "E18" "  -2.67%  "

# Row 19
Set-TextCell $ws "B19" "WrappedBTC"
Set-TextCell $ws "C19" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws "D19" "69.466.63"
Set-TextCell $ws "E19" "  -1.81%  "

# Row 20
Set-TextCell $ws "E20" "  +1.22%  "

# Row 21
Set-TextCell $ws "D21" "17.39"
Set-TextCell $ws "E21" "  -2.04%  "

# Row 22
Set-TextCell $ws "D22" "11.20"
Set-TextCell $ws "E22" "  -1.78%  "

# Row 23
Set-TextCell $ws "D23" "0.891"
Set-TextCell $ws "E23" "  -4.03%  "

# Row 24
Set-TextCell $ws "D24" "15.97"
Set-TextCell $ws "E24" "  -6.62%  "

# Row 25
Set-TextCell $ws "D25" "97.31"
Set-TextCell $ws "E25" "  -3.08%  "

# Row 26
Set-TextCell $ws "D26" "3.80"
Set-TextCell $ws "E26" "  -2.54%  "

# Row 27
Set-TextCell $ws "E27" "  +0.86%  "

# Row 28
Set-TextCell $ws "D28" "0.999"
Set-TextCell $ws "E28" "  +0.01%  "

# Row 29
Set-TextCell $ws "D29" "2.64"
Set-TextCell $ws "E29" "  -5.25%  "

# Row 30
Set-TextCell $ws "D30" "9.33"
Set-TextCell $ws "E30" "  -6.52%  "

# Row 31
Set-TextCell $ws "D31" "32.89"
Set-TextCell $ws "E31" "  -5.44%  "

# Row 32
Set-TextCell $ws "E32" "  -5.87%  "

# Row 33
Set-TextCell $ws "D33" "8.57"
Set-TextCell $ws "E33" "  -4.76%  "

# Row 34
Set-TextCell $ws "D34" "1.34"
Set-TextCell $ws "E34" "  -3.92%  "

# Row 35
Set-TextCell $ws "D35" "7.00"
Set-TextCell $ws "E35" "  -3.99%  "

# Row 36
Set-TextCell $ws "D36" "639.02"
Set-TextCell $ws "E36" "  +9.64%  "

# Row 37
Set-TextCell $ws "D37" "10.80"
Set-TextCell $ws "E37" "  -2.20%  "

# Row 38
Set-TextCell $ws "D38" "3.51"
Set-TextCell $ws "E38" "  -12.15%  "

# Row 39
Set-TextCell $ws "E39" "  -3.80%  "

# Row 40
Set-TextCell $ws "D40" "57.32"
Set-TextCell $ws "E40" "  -1.68%  "

# Row 41
Set-TextCell $ws "E41" "  +0.05%  "

# Row 42
Set-TextCell $ws "D42" "0.0457"
Set-TextCell $ws "E42" "  +0.53%  "

# Row 43
Set-TextCell $ws "D43" "0.136"
Set-TextCell $ws "E43" "  -3.37%  "

# Row 44
Set-TextCell $ws "D44" "3.395.26"
Set-TextCell $ws "E44" "  -5.26%  "

# Row 45
Set-TextCell $ws "D45" "0.329"
Set-TextCell $ws "E45" "  -4.22%  "

# Row 46
Set-TextCell $ws "D46" "0.0₃0703"
Set-TextCell $ws "E46" "  -5.44%  "

# Row 47
Set-TextCell $ws "D47" "32.79"
Set-TextCell $ws "E47" "  -5.55%  "

# Row 48
Set-TextCell $ws "D48" "2.56"
Set-TextCell $ws "E48" "  -5.33%  "

# Row 49
Set-TextCell $ws "D49" "2.75"
Set-TextCell $ws "E49" "  -4.24%  "

# Row 50
Set-TextCell $ws "E50" "  -2.17%  "

# Row 51
Set-TextCell $ws "D51" "132.66"
Set-TextCell $ws "E51" "  -1.70%  "
